$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2 = 62.47613131263947
    3 = 66.53155233183062
    4 = 68.77501848443904
    5 = 69.39889936748489
    6 = 75.54199288655842
    7 = 74.5336110899503
    8 = 71.88476462614811
}

foreach ($row in $values.Keys) {
    $ws.Range("B$row").Value = $values[$row]
    $ws.Range("C$row").Formula = "=AVERAGE(B$row`:B$row)"
    $ws.Range("D$row").Formula = "=STDEV(B$row`:B$row)"
}

$wb.Save()
